# Auto-generated from the unified diff: updates crypto price/volume cells
# and restores the two swapped coin-row pairs to their new order/values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.431.51"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "'1.790.97"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'226.14"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'0.557"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'32.62"
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("D9").Value = "'0.297"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("D10").Value = "'0.0691"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'2.049.53"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.798.48"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.08"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "'34.404.88"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "'68.77"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "'246.80"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "'0.0₃0797"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "'164.48"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'16.49"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.80"
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("E33").Value = "  +6.93%  "
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("D35").Value = "'1.419.16"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("E36").Value = "  +5.08%  "
$ws.Range("D37").Value = "'0.669"
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0192"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +5.23%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("D44").Value = "'13.63"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("D46").Value = "'6.06"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "'1.946.49"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₆0131"
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.02%  "
